# Adds two new columns, I ("I0") and J ("IF"), to the sheet with header
# style matching the existing headers (bold, centered, thin-bordered),
# and fills in the per-row numeric values for rows 2-83.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells -----------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold/centered/bordered) from the existing H1
# header cell onto the two new header cells so they share the same
# cell style as the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data values for columns I (I0) and J (IF), rows 2-83 -------------
$i0Values = @(8,8,8,7,8,9,9,8,7,8,8,7,5,7,7,7,7,7,6,9,8,8,7,7,8,8,6,7,7,6,8,7,6,9,8,6,6,8,7,11,10,6,6,6,7,10,7,9,8,7,7,7,7,4,6,6,6,5,8,6,6,4,9,6,7,4,4,6,8,6,9,7,5,9,7,4,7,7,6,3,7,3)
$ifValues = @(8,8,8,7,8,9,9,8,8,8,8,8,6,8,8,8,7,7,7,9,9,8,7,7,8,8,6,7,8,6,8,7,7,9,8,7,7,8,7,12,10,7,6,6,7,10,7,9,8,7,7,7,7,6,7,6,6,6,8,6,6,4,9,6,8,4,5,6,8,6,9,8,6,9,7,4,7,7,6,4,7,3)

$startRow = 2
for ($idx = 0; $idx -lt $i0Values.Length; $idx++) {
    $row = $startRow + $idx
    $ws.Cells.Item($row, 9).Value = $i0Values[$idx]
    $ws.Cells.Item($row, 10).Value = $ifValues[$idx]
}
